# Automatic update of files.
# Update the "Förändrad" (Changed) date column (C) for rows 2-6
# from serial date 45233 (2023-11-03) to 45243 (2023-11-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 3).Value = 45243
}
